# Updates cryptos list prices/volumes per GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (matches the source data, which is all typed as inline/shared strings,
    # even for cells that look numeric, e.g. "1.00" or "34.157.65").
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "34.157.65"
Set-TextCell "E2" "  +0.77%  "
Set-TextCell "D3" "1.787.70"
Set-TextCell "E3" "  -1.22%  "
Set-TextCell "E4" "  +0.03%  "
Set-TextCell "D5" "226.36"
Set-TextCell "E5" "  -0.63%  "
Set-TextCell "D6" "0.551"
Set-TextCell "E6" "  +1.75%  "
Set-TextCell "E7" "  +0.01%  "
Set-TextCell "D8" "31.56"
Set-TextCell "E8" "  +2.19%  "
Set-TextCell "B9" "Cardano"
Set-TextCell "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell "D9" "0.282"
Set-TextCell "E9" "  +1.26%  "
Set-TextCell "B10" "Dogecoin"
Set-TextCell "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D10" "0.0661"
Set-TextCell "E10" "  -0.34%  "
Set-TextCell "B11" "TRON"
Set-TextCell "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D11" "0.0931"
Set-TextCell "E11" "  +0.10%  "
Set-TextCell "B12" "WrappedliquidstakedEther2.0"
Set-TextCell "C12" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D12" "2.045.56"
Set-TextCell "E12" "  -1.22%  "
Set-TextCell "B13" "Chainlink"
Set-TextCell "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D13" "11.32"
Set-TextCell "E13" "  +12.60%  "
Set-TextCell "B14" "WrappedEther"
Set-TextCell "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D14" "1.784.03"
Set-TextCell "E14" "  -1.46%  "
Set-TextCell "B15" "Polygon"
Set-TextCell "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D15" "0.630"
Set-TextCell "E15" "  -0.68%  "
Set-TextCell "B16" "WrappedBTC"
Set-TextCell "C16" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D16" "34.151.25"
Set-TextCell "E16" "  +0.84%  "
Set-TextCell "B17" "Polkadot"
Set-TextCell "C17" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D17" "4.24"
Set-TextCell "E17" "  +0.16%  "
Set-TextCell "B18" "Litecoin"
Set-TextCell "C18" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D18" "69.39"
Set-TextCell "E18" "  +0.45%  "
Set-TextCell "B19" "BitcoinCash"
Set-TextCell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D19" "254.59"
Set-TextCell "E19" "  -0.01%  "
Set-TextCell "B20" "ShibaInu"
Set-TextCell "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D20" "0.0₃0742"
Set-TextCell "E20" "  +0.28%  "
Set-TextCell "B21" "Dai"
Set-TextCell "C21" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D21" "1.00"
Set-TextCell "E21" "  -0.02%  "
Set-TextCell "B22" "Avalanche"
Set-TextCell "C22" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D22" "10.50"
Set-TextCell "E22" "  +1.37%  "
Set-TextCell "B23" "Uniswap"
Set-TextCell "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D23" "4.22"
Set-TextCell "E23" "  -1.45%  "
Set-TextCell "B24" "Toncoin"
Set-TextCell "C24" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D24" "2.14"
Set-TextCell "E24" "  -1.75%  "
Set-TextCell "B25" "Monero"
Set-TextCell "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D25" "156.30"
Set-TextCell "E25" "  -1.56%  "
Set-TextCell "B26" "EthereumClassic"
Set-TextCell "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D26" "16.53"
Set-TextCell "E26" "  +0.59%  "
Set-TextCell "B27" "Cosmos"
Set-TextCell "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D27" "7.03"
Set-TextCell "E27" "  +0.18%  "
Set-TextCell "B28" "Stellar"
Set-TextCell "C28" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D28" "0.114"
Set-TextCell "E28" "  -0.18%  "
Set-TextCell "B29" "BinanceUSD"
Set-TextCell "C29" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D29" "1.00"
Set-TextCell "E29" "  +0.05%  "
Set-TextCell "B30" "Filecoin"
Set-TextCell "C30" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D30" "3.80"
Set-TextCell "E30" "  +0.10%  "
Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.0518"
Set-TextCell "E31" "  +2.15%  "
Set-TextCell "B32" "PancakeSwap"
Set-TextCell "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D32" "1.20"
Set-TextCell "E32" "  +0.23%  "
Set-TextCell "B33" "InternetComputer(DFINITY)"
Set-TextCell "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D33" "3.58"
Set-TextCell "E33" "  +2.22%  "
Set-TextCell "B34" "LidoDAOToken"
Set-TextCell "C34" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D34" "1.84"
Set-TextCell "E34" "  +2.53%  "
Set-TextCell "B35" "Maker"
Set-TextCell "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D35" "1.452.17"
Set-TextCell "E35" "  -5.89%  "
Set-TextCell "B36" "TrustWalletToken"
Set-TextCell "C36" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D36" "1.07"
Set-TextCell "E36" "  -0.33%  "
Set-TextCell "B37" "ImmutableX"
Set-TextCell "C37" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D37" "0.634"
Set-TextCell "E37" "  +3.00%  "
Set-TextCell "B38" "VeChain"
Set-TextCell "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D38" "0.0188"
Set-TextCell "E38" "  +1.16%  "
Set-TextCell "B39" "Aave"
Set-TextCell "C39" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D39" "83.44"
Set-TextCell "E39" "  -0.06%  "
Set-TextCell "B40" "MXToken"
Set-TextCell "C40" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D40" "2.86"
Set-TextCell "E40" "  +1.15%  "
Set-TextCell "B41" "HuobiToken"
Set-TextCell "C41" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D41" "2.35"
Set-TextCell "E41" "  +0.39%  "
Set-TextCell "B42" "ARBITRUM"
Set-TextCell "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D42" "0.898"
Set-TextCell "E42" "  -0.14%  "
Set-TextCell "B43" "RenderToken"
Set-TextCell "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D43" "2.07"
Set-TextCell "E43" "  -1.13%  "
Set-TextCell "B44" "Kaspa"
Set-TextCell "C44" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D44" "0.0509"
Set-TextCell "E44" "  -2.43%  "
Set-TextCell "B45" "WEMIXToken"
Set-TextCell "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D45" "1.07"
Set-TextCell "E45" "  -0.87%  "
Set-TextCell "B46" "FraxShare"
Set-TextCell "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D46" "5.87"
Set-TextCell "E46" "  +4.14%  "
Set-TextCell "B47" "RocketPoolETH"
Set-TextCell "C47" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D47" "1.944.08"
Set-TextCell "E47" "  -0.95%  "
Set-TextCell "B48" "PaxDollar"
Set-TextCell "C48" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D48" "1.00"
Set-TextCell "E48" "  +0.06%  "
Set-TextCell "D49" "11.92"
Set-TextCell "E49" "  +7.19%  "
Set-TextCell "B50" "BitcoinSV"
Set-TextCell "C50" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell "D50" "50.51"
Set-TextCell "E50" "  -3.04%  "
Set-TextCell "B51" "Quant"
Set-TextCell "C51" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D51" "98.06"
Set-TextCell "E51" "  +2.48%  "
